# Update the arithmetic problems in the single table on the page.
# The table is 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17)
# holds problem text, the others are spacer rows. Addressing cells
# positionally (Table.Cell(row, col)) avoids any ambiguity from the
# duplicate "65÷9=" text that appears twice in row 9 but must map to two
# different replacement values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "12÷9="; New = "85÷2=" },
    @{ Row = 1;  Col = 2; Old = "10÷7="; New = "87÷4=" },
    @{ Row = 1;  Col = 3; Old = "19÷6="; New = "55÷4=" },
    @{ Row = 1;  Col = 4; Old = "70÷2="; New = "39÷3=" },
    @{ Row = 1;  Col = 5; Old = "18÷5="; New = "49÷9=" },

    @{ Row = 5;  Col = 1; Old = "40÷5="; New = "64÷3=" },
    @{ Row = 5;  Col = 2; Old = "66÷9="; New = "10÷8=" },
    @{ Row = 5;  Col = 3; Old = "65÷7="; New = "53÷8=" },
    @{ Row = 5;  Col = 4; Old = "56÷3="; New = "99÷9=" },
    @{ Row = 5;  Col = 5; Old = "95÷6="; New = "28÷8=" },

    @{ Row = 9;  Col = 1; Old = "65÷9="; New = "33÷3=" },
    @{ Row = 9;  Col = 2; Old = "70÷7="; New = "73÷3=" },
    @{ Row = 9;  Col = 3; Old = "85÷7="; New = "30÷6=" },
    @{ Row = 9;  Col = 4; Old = "22÷9="; New = "23÷9=" },
    @{ Row = 9;  Col = 5; Old = "65÷9="; New = "23÷9=" },

    @{ Row = 13; Col = 1; Old = "50÷7="; New = "80÷3=" },
    @{ Row = 13; Col = 2; Old = "71÷9="; New = "40÷5=" },
    @{ Row = 13; Col = 3; Old = "86÷9="; New = "13÷5=" },
    @{ Row = 13; Col = 4; Old = "82÷4="; New = "61÷4=" },
    @{ Row = 13; Col = 5; Old = "43÷4="; New = "66÷4=" },

    @{ Row = 17; Col = 1; Old = "21÷7="; New = "47÷2=" },
    @{ Row = 17; Col = 2; Old = "24÷8="; New = "82÷2=" },
    @{ Row = 17; Col = 3; Old = "80÷6="; New = "86÷6=" },
    @{ Row = 17; Col = 4; Old = "69÷9="; New = "56÷4=" },
    @{ Row = 17; Col = 5; Old = "90÷6="; New = "58÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing end-of-cell marker Word appends to Range.Text.
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}
